$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-04-13 Sunday", $false, $false, $false, $false,
                         $false, $true, 1, $false, "2025-04-14 Monday", 2)

# Update the division problems in the table (5 rows x 5 columns of values)
$t = $d.Tables.Item(1)

$values = @(
    @("47÷4=", "31÷2=", "42÷4=", "94÷9=", "36÷7="),
    @("70÷8=", "61÷2=", "85÷3=", "31÷5=", "20÷2="),
    @("56÷2=", "25÷6=", "69÷7=", "47÷2=", "24÷3="),
    @("84÷3=", "63÷2=", "57÷4=", "79÷8=", "94÷9="),
    @("34÷9=", "15÷9=", "50÷8=", "53÷8=", "93÷8=")
)

$rowIndexes = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $rowIndexes.Length; $i++) {
    $rowIdx = $rowIndexes[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($rowIdx, $c)
        $cell.Range.Text = $values[$i][$c - 1]
    }
}
